# Daily auto-update routine.
#
# The sheet tracks, per shop (row), a rental/service period:
#   D = total days, E = remaining days, F = start date (yyyyMMdd, plain number).
# Each day the sheet is refreshed against "today": E is recomputed as the
# number of days left until (F + D). When a period has lapsed (remaining
# would hit zero or go negative) the period is treated as renewed starting
# "today", so F is reset to today and E is reset back to the full D.
#
# This run advances the reference day by one (2026-02-02 -> 2026-02-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function ParseYmd($ymdVal) {
    $s = [string]$ymdVal
    $y = [int]$s.Substring(0, 4)
    $m = [int]$s.Substring(4, 2)
    $d = [int]$s.Substring(6, 2)
    $dt = Get-Date -Year $y -Month $m -Day $d
    return $dt
}

$newToday = Get-Date -Year 2026 -Month 2 -Day 3
$newTodayYmd = [int]$newToday.ToString("yyyyMMdd")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $dVal = $dCell.Value()
    $eVal = $eCell.Value()
    $fVal = $fCell.Value()

    if ($dVal -eq $null) { continue }
    if ($eVal -eq $null) { continue }
    if ($fVal -eq $null) { continue }

    $fStr = [string][int64]$fVal
    if ($fStr.Length -ne 8) {
        # malformed start date (data-entry typo) - leave untouched
        continue
    }

    $start = ParseYmd($fVal)
    $end = $start.AddDays([double]$dVal)

    $newRemaining = [int][math]::Round($end.ToOADate() - $newToday.ToOADate())

    if ($newRemaining -le 0) {
        $eCell.Value = [int]$dVal
        $fCell.Value = $newTodayYmd
    } else {
        $eCell.Value = $newRemaining
    }
}
